$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update date strings in column A (rows 3-21): change "/" separators to "-"
$dates = @{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

foreach ($row in $dates.Keys) {
    $cell = $ws.Cells.Item($row, 1)
    # Force text interpretation so Excel doesn't auto-convert ambiguous
    # "DD-MM-YYYY" strings (day <= 12) into real date serials.
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$row]
    # Restore the cell to its original (unstyled) Normal style so no
    # spurious style/format change is introduced.
    $cell.Style = "Normal"
}

# Update attendance counts for row 14 (05-09-2022): Real/Total attendance recorded, Absent cleared
$ws.Cells.Item(14, 4).Value = 1
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 8).Value = 0

# Update attendance counts for row 20 (26-09-2022): Real/Total attendance recorded, Absent cleared
$ws.Cells.Item(20, 4).Value = 1
$ws.Cells.Item(20, 5).Value = 1
$ws.Cells.Item(20, 8).Value = 0
